$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.451.41"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.865.29"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'235.15"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4822"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.2797"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").Value = "'0.06498"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "1.947.76"
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("D11").Value = "'0.07431"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").Value = "'16.31"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'5.064"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "'87.08"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'0.6454"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "30.408.50"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "'12.97"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'233.10"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007539"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "2.113.42"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.154"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "'6.093"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'9.324"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'166.89"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "'18.33"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").Value = "'1.920"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "'0.1021"
$ws.Range("E29").Value = "  +9.23%  "
$ws.Range("D30").Value = "'1.373"
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("D31").Value = "'4.264"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "'3.988"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'0.04969"
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "'1.179"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").Value = "'0.7301"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'2.717"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'0.01911"
$ws.Range("E38").Value = "  +4.16%  "
$ws.Range("D39").Value = "'2.629"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "'0.9107"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "'2.037"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").Value = "'105.98"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'0.9960"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "'0.4193"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'5.548"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("D46").Value = "'7.214"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "'61.98"
$ws.Range("E47").Value = "  -6.30%  "
$ws.Range("D48").Value = "'0.1231"
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").Value = "'8.826"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "'1.438"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").Value = "'33.53"
$ws.Range("E51").Value = "  -2.36%  "
